# Applies the changes described by the diff:
#  1) Update the DevExpress evaluation-warning string version number
#     from v25.1.7.0 to v25.2.3.0 (on the "Evaluation Warning" sheet).
#  2) Update the data rows on the "Sheet" worksheet:
#       Row 2: A2 55 -> 482, I2 44 -> 252, J2 46009.370355262974 -> 46025.683620825424
#       Row 3: A3 56 -> 483, I3 45 -> 253, J3 46009.370355262974 -> 46025.683620825424

$wb = $excel.ActiveWorkbook

# --- 1. Fix the "Evaluation Warning" text ---
$wsWarning = $wb.Worksheets.Item("Evaluation Warning")
$wsWarning.Range("A6").Value = "or purchase a new license (devexpress.com/BUY) to continue use of DevExpress product libraries (v25.2.3.0)."

# --- 2. Update the data table on the "Sheet" worksheet ---
$ws = $wb.Worksheets.Item("Sheet")

$ws.Range("A2").Value = 482
$ws.Range("I2").Value = 252
$ws.Range("J2").Value = 46025.683620825424

$ws.Range("A3").Value = 483
$ws.Range("I3").Value = 253
$ws.Range("J3").Value = 46025.683620825424
